$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.744.22"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "3.486.80"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.431"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").Value = "4.091.46"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").Value = "66.787.18"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").Value = "3.507.09"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "392.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.535"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("E29").Value = "  -2.94%  "
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.00%  "
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0739"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("D43").Value = "2.798.64"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("E45").Value = "  +1.90%  "
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "335.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.28%  "
